# Update "Fruta, Vega Modelo de Temuco - Níspero" weekly data.
# The underlying data rows (2-20) were re-sorted / refreshed with new
# period values; columns D (Fecha), M (Volumen), N (Precio mínimo),
# O (Precio máximo), P (Precio promedio ponderado),
# Q (Unidad de comercialización), R (Origen), S (Precio $/Kg) and
# T (Kg / unidad) change per-row as described below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number -> ordered values for columns D,M,N,O,P,Q,R,S,T
$rows = @{
    2  = @(44511, 45,  28000, 28000, 28000, '$/bandeja 10 kilos', 'Provincia de Los Andes', 2800, 10)
    3  = @(44511, 45,  3200,  3200,  3200,  '$/bandeja 10 kilos', 'Provincia de Quillota',  320,  10)
    5  = @(44519, 30,  28000, 28000, 28000, '$/bandeja 10 kilos', 'Provincia de Quillota',  2800, 10)
    6  = @(44515, 80,  28000, 28000, 28000, '$/bandeja 10 kilos', 'Provincia de Los Andes', 2800, 10)
    7  = @(44503, 50,  28000, 28000, 28000, '$/bandeja 10 kilos', 'Provincia de Quillota',  2800, 10)
    8  = @(44496, 55,  28000, 28000, 28000, '$/bandeja 10 kilos', 'Provincia de Quillota',  2800, 10)
    10 = @(44868, 30,  14000, 14000, 14000, '$/bandeja 5 kilos',  'Provincia de Quillota',  2800, 5)
    11 = @(44858, 90,  20000, 20000, 20000, '$/bandeja 5 kilos',  'Provincia de Quillota',  4000, 5)
    12 = @(44488, 100, 12000, 12000, 12000, '$/bandeja 5 kilos',  'La Ligua',               2400, 5)
    13 = @(44879, 25,  30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Quillota',  3000, 10)
    14 = @(44483, 35,  10000, 10000, 10000, '$/bandeja 5 kilos',  'Provincia de Quillota',  2000, 5)
    15 = @(44466, 80,  11000, 11000, 11000, '$/bandeja 5 kilos',  'La Ligua',               2200, 5)
    16 = @(44902, 90,  25000, 25000, 25000, '$/bandeja 10 kilos', 'Provincia de Quillota',  2500, 10)
    17 = @(44921, 55,  15000, 15000, 15000, '$/bandeja 7 kilos',  'Provincia de Quillota',  2143, 7)
    18 = @(44166, 20,  12000, 12000, 12000, '$/caja 18 kilos',    'La Ligua',               667,  18)
    19 = @(44901, 40,  25000, 25000, 25000, '$/bandeja 10 kilos', 'Provincia de Quillota',  2500, 10)
    20 = @(44889, 50,  30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Quillota',  3000, 10)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("M$r").Value = $vals[1]
    $ws.Range("N$r").Value = $vals[2]
    $ws.Range("O$r").Value = $vals[3]
    $ws.Range("P$r").Value = $vals[4]
    $ws.Range("Q$r").Value = $vals[5]
    $ws.Range("R$r").Value = $vals[6]
    $ws.Range("S$r").Value = $vals[7]
    $ws.Range("T$r").Value = $vals[8]
}
